# Add new conformance files (C042 "imir" test, multilayer005, and the
# matching B025 bitstream row) to the "Test Content" conformance table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert a new row for "C042" right before "multilayer001" (old row 42)
# ---------------------------------------------------------------------
$ws.Rows.Item(42).Insert()
$ws.Range("A42").Value2 = "C042"
$ws.Range("B42").Value2 = 'An image item vertically mirrored with an "imir" property.'
$ws.Range("C42").Value2 = "B001"
$ws.Range("D42").Value2 = "heic, mif1"
$ws.Range("B42").WrapText = $true

# ---------------------------------------------------------------------
# 2) Insert a new row for "multilayer005" right after "multilayer004"
#    (which, after the previous insert, now lives at row 46)
# ---------------------------------------------------------------------
$ws.Rows.Item(47).Insert()
$ws.Range("A47").Value2 = "multilayer005"
$ws.Range("B47").Value2 = "A multi-layer multiview file with 'ster' grouping."
$ws.Range("C47").Value2 = "B025"
$ws.Range("D47").Value2 = "heic, heis, mif1"
$ws.Range("B47").WrapText = $true

# ---------------------------------------------------------------------
# 3) Append the matching bitstream description row "B025" at the very
#    end of the "Bitstream ID" table (new last row, 76)
# ---------------------------------------------------------------------
$ws.Range("A76").Value2 = "B025"
$ws.Range("B76").Value2 = "1 HEVC encoded multi-layer frame. Layer 0 is the left view, layer 1 is the right view.  (512x256 resolution)"
$ws.Range("C76").Value2 = "B22"
$ws.Range("B76").WrapText = $true

# Leave the selection where the author's session ended up.
$ws.Range("B77").Select() | Out-Null
